$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.07443033333333333
$ws.Range("H2").Value = 0.223291
$ws.Range("I2").Value = 0.07586947613633815
$ws.Range("J2").Value = 0.07586947613633817
$ws.Range("M2").Value = 0.1801153333333333
$ws.Range("N2").Value = 0.540346
$ws.Range("O2").Value = 0.01663333613045927
$ws.Range("P2").Value = 0.01663333613045927
$ws.Range("Q2").Value = 0.01340604429844444
$ws.Range("R2").Value = 0.120654398686
$ws.Range("S2").Value = 0.001261962498617571
$ws.Range("T2").Value = 0.001261962498617571
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.07443033333333333
$ws.Range("H3").Value = 0.223291
$ws.Range("I3").Value = 0.07586947613633815
$ws.Range("J3").Value = 0.07586947613633817
$ws.Range("O3").Value = 0.407089716880577
$ws.Range("P3").Value = 0.407089716880577
$ws.Range("Q3").Value = 0.3281039194505556
$ws.Range("R3").Value = 2.952935275055
$ws.Range("S3").Value = 0.0308856835602196
$ws.Range("T3").Value = 0.0308856835602196
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.07443033333333333
$ws.Range("H4").Value = 0.223291
$ws.Range("I4").Value = 0.07586947613633815
$ws.Range("J4").Value = 0.07586947613633817
$ws.Range("M4").Value = 6.240258333333333
$ws.Range("N4").Value = 18.720775
$ws.Range("O4").Value = 0.5762769469889637
$ws.Range("P4").Value = 0.5762769469889637
$ws.Range("Q4").Value = 0.4644645078361111
$ws.Range("R4").Value = 4.180180570525001
$ws.Range("S4").Value = 0.04372183007750099
$ws.Range("T4").Value = 0.04372183007750099
$ws.Range("I5").Value = 0.7501574873245639
$ws.Range("J5").Value = 0.7501574873245638
$ws.Range("M5").Value = 0.1801153333333333
$ws.Range("N5").Value = 0.540346
$ws.Range("O5").Value = 0.01663333613045927
$ws.Range("P5").Value = 0.01663333613045927
$ws.Range("Q5").Value = 0.1325519170293333
$ws.Range("R5").Value = 1.192967253264
$ws.Range("S5").Value = 0.01247762163745021
$ws.Range("T5").Value = 0.01247762163745021
$ws.Range("I6").Value = 0.7501574873245639
$ws.Range("J6").Value = 0.7501574873245638
$ws.Range("O6").Value = 0.407089716880577
$ws.Range("P6").Value = 0.407089716880577
$ws.Range("S6").Value = 0.3053813991308018
$ws.Range("T6").Value = 0.3053813991308018
$ws.Range("I7").Value = 0.7501574873245639
$ws.Range("J7").Value = 0.7501574873245638
$ws.Range("M7").Value = 6.240258333333333
$ws.Range("N7").Value = 18.720775
$ws.Range("O7").Value = 0.5762769469889637
$ws.Range("P7").Value = 0.5762769469889637
$ws.Range("Q7").Value = 4.592380834733333
$ws.Range("R7").Value = 41.33142751259999
$ws.Range("S7").Value = 0.4322984665563119
$ws.Range("T7").Value = 0.4322984665563118
$ws.Range("G8").Value = 0.170673
$ws.Range("H8").Value = 0.512019
$ws.Range("I8").Value = 0.173973036539098
$ws.Range("J8").Value = 0.173973036539098
$ws.Range("M8").Value = 0.1801153333333333
$ws.Range("N8").Value = 0.540346
$ws.Range("O8").Value = 0.01663333613045927
$ws.Range("P8").Value = 0.01663333613045927
$ws.Range("Q8").Value = 0.030740824286
$ws.Range("R8").Value = 0.276667418574
$ws.Range("S8").Value = 0.002893751994391489
$ws.Range("T8").Value = 0.002893751994391489
$ws.Range("G9").Value = 0.170673
$ws.Range("H9").Value = 0.512019
$ws.Range("I9").Value = 0.173973036539098
$ws.Range("J9").Value = 0.173973036539098
$ws.Range("O9").Value = 0.407089716880577
$ws.Range("P9").Value = 0.407089716880577
$ws.Range("Q9").Value = 0.752361003055
$ws.Range("R9").Value = 6.771249027495
$ws.Range("S9").Value = 0.07082263418955569
$ws.Range("T9").Value = 0.07082263418955569
$ws.Range("G10").Value = 0.170673
$ws.Range("H10").Value = 0.512019
$ws.Range("I10").Value = 0.173973036539098
$ws.Range("J10").Value = 0.173973036539098
$ws.Range("M10").Value = 6.240258333333333
$ws.Range("N10").Value = 18.720775
$ws.Range("O10").Value = 0.5762769469889637
$ws.Range("P10").Value = 0.5762769469889637
$ws.Range("Q10").Value = 1.065043610525
$ws.Range("R10").Value = 9.585392494724999
$ws.Range("S10").Value = 0.1002566503551508
$ws.Range("T10").Value = 0.1002566503551508
